$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 80: Nov 4, 2021 - HOLIDAY
$ws.Range("A80").Value = "04/11/2021"
$ws.Range("B80").Value = "HOLIDAY"

# Row 81: Nov 5, 2021 - Status
$ws.Range("A81").Value = "05/11/2021"
$ws.Range("B81").Value = "Continued on codec 2.0 "
$ws.Range("C81").Value = "Revise the studied concepts and improve my understanding"
$ws.Range("D81").Value = "Listed the doubts , have to elaborate "

# Row 82
$ws.Range("B82").Value = "Not a considerable progress but OPENMAX is reviewed and got understanding on the IL"
$ws.Range("D82").Value = "Revision of C-DS-OS concepts"

# Row 83
$ws.Range("B83").Value = "Listed and observed the OPENMAX APIs"

# Row 84
$ws.Range("B84").Value = "Following android media architechture , observed abstractions in OPENMAX "

$ws.Range("D82").Select() | Out-Null
